$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Feuil1" (sheet2): reorder / extend the second lookup list (rows 11-19)
# ---------------------------------------------------------------------------
$f1 = $wb.Worksheets.Item("Feuil1")

# D11 used to hold "RT.ART" - that entry now lives in the LIST sheet list,
# so clear it from here.
$f1.Range("D11").ClearContents()

$f1.Range("B14").Value = "AD.DEP.001.FON.01"
$f1.Range("B15").Value = "RO.ACT"
$f1.Range("B16").Value = "RO.FOU"

$f1.Range("B17").Value = "MP.CPT"
$f1.Range("B17").NumberFormat = "@"

$f1.Range("B18").Value = "RT.ART"
$f1.Range("B18").NumberFormat = "@"

$f1.Range("B19").Value = "AD.SEC.014.FON.01"

$f1.Range("D11").Select()

# ---------------------------------------------------------------------------
# Sheet "LIST" (sheet1): rebuild the column A list (rows 2-10)
# ---------------------------------------------------------------------------
$list = $wb.Worksheets.Item("LIST")

$list.Range("A2").Value = "AD.SEC.001.FON.02"
$list.Range("A2").NumberFormat = "@"

$list.Range("A3").Value = "AD.SEC.001.FON.01"
$list.Range("A3").NumberFormat = "@"

$list.Range("A4").Value = "AD.SEC.001.FON.03"
$list.Range("A4").NumberFormat = "@"

$list.Range("A5").Value = "AD.DEP.001.FON.01"
$list.Range("A5").NumberFormat = "@"

$list.Range("A6").Value = "RO.ACT"
$list.Range("A6").NumberFormat = "@"

$list.Range("A7").Value = "RO.FOU"
$list.Range("A7").NumberFormat = "@"

$list.Range("A8").Value = "MP.CPT"
$list.Range("A8").NumberFormat = "@"

$list.Range("A9").Value = "RT.ART"
$list.Range("A9").NumberFormat = "@"

$list.Range("A10").Value = "AD.SEC.014.FON.01"

$list.Range("A2:A10").Select()
